$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 0.606
$ws.Range("E2").Value = 0.805
$ws.Range("F2").Value = 0.798
$ws.Range("G2").Value = 0.369
$ws.Range("H2").Value = 1.528

# Row 3
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 0.606
$ws.Range("E3").Value = 0.966
$ws.Range("F3").Value = 0.893
$ws.Range("G3").Value = 0.22
$ws.Range("H3").Value = 1.708

# Row 4
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 0.613
$ws.Range("E4").Value = 0.984
$ws.Range("F4").Value = 1.037
$ws.Range("G4").Value = 0.256
$ws.Range("H4").Value = 1.609

# Row 5
$ws.Range("D5").Value = 0.762
$ws.Range("E5").Value = 0.75
$ws.Range("F5").Value = 0.75
$ws.Range("G5").Value = 0.47
$ws.Range("H5").Value = 2.074

# Row 6
$ws.Range("D6").Value = 0.618
$ws.Range("E6").Value = 1.141
$ws.Range("F6").Value = 1.328
$ws.Range("G6").Value = 0.257
$ws.Range("H6").Value = 1.949
